# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted at row 141 (pushing the existing
# rows 141-195 down to 142-196). The new record is for "Primera" quality
# Sandia, reported in $/kilo (volumen en unidades) from Peru.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 141, shifting the remaining
# data (old rows 141-195) down to rows 142-196.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A141").Value = 9
$ws.Range("B141").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C141").Value = "Metropolitana"
$ws.Range("D141").Value = 44452
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 100112028
$ws.Range("G141").Value = "Sandia"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 160
$ws.Range("K141").Value = 1500
$ws.Range("L141").Value = 1500
$ws.Range("M141").Value = 1500
$ws.Range("N141").Value = "$/kilo (volumen en unidades)"
$ws.Range("O141").Value = "Perú"
$ws.Range("P141").Value = 1500
$ws.Range("Q141").Value = 1
$ws.Range("R141").Value = "Hortaliza"
